# Refresh the cryptos list: update each coin's Price (col D) and
# Volume(1h) (col E) cell to the latest scraped text.
#
# Both columns hold plain text in the source data (e.g. "27.666.21",
# "0.0611", "  +1.03%  "), not numbers. Excel's Range.Value setter
# auto-detects values that *look* like a genuine number (no thousands-style
# dots, a single decimal point) and silently coerces them to a numeric
# cell, which would corrupt values like "0.524" (-> 0.524 as a float,
# losing the literal text) and flips the cell's type. To keep those cells
# as text -- exactly like the source -- we flip NumberFormat to "@" (Text)
# immediately before assigning, then clear the formatting again right
# after so the cell keeps the workbook's default (un-styled) appearance.
# Values that already fail Excel's numeric parse (e.g. "27.666.21", which
# has two dots) or are not number-like (the "  +x.xx%  " strings) don't
# need this and are just assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" '27.666.21'
$ws.Range("E2").Value = '  +1.03%  '

# Row 3 - Ethereum
Set-TextValue "D3" '1.639.46'
$ws.Range("E3").Value = '  +0.06%  '

# Row 5 - BNB
Set-TextValue "D5" '212.83'
$ws.Range("E5").Value = '  +0.61%  '

# Row 6 - XRP
Set-TextValue "D6" '0.524'
$ws.Range("E6").Value = '  -0.33%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  +0.04%  '

# Row 8 - Solana
$ws.Range("E8").Value = '  +1.26%  '

# Row 9 - Cardano
$ws.Range("E9").Value = '  +0.76%  '

# Row 10 - Dogecoin
Set-TextValue "D10" '0.0611'
$ws.Range("E10").Value = '  +0.17%  '

# Row 11 - TRON
Set-TextValue "D11" '0.0895'
$ws.Range("E11").Value = '  +0.32%  '

# Row 12 - Wrapped liquid staked Ether 2.0
Set-TextValue "D12" '1.872.00'
$ws.Range("E12").Value = '  +0.09%  '

# Row 13 - Wrapped Ether
Set-TextValue "D13" '1.625.95'
$ws.Range("E13").Value = '  -0.62%  '

# Row 14 - Polkadot
$ws.Range("E14").Value = '  +0.52%  '

# Row 15 - Polygon
$ws.Range("E15").Value = '  +0.01%  '

# Row 16 - Litecoin
Set-TextValue "D16" '64.65'
$ws.Range("E16").Value = '  +0.59%  '

# Row 17 - Wrapped BTC
Set-TextValue "D17" '27.644.39'

# Row 18 - Bitcoin Cash
Set-TextValue "D18" '230.30'
$ws.Range("E18").Value = '  +0.74%  '

# Row 19 - Chainlink
Set-TextValue "D19" '7.70'
$ws.Range("E19").Value = '  +2.30%  '

# Row 20 - Shiba Inu
$ws.Range("E20").Value = '  +0.53%  '

# Row 21 - Dai
$ws.Range("E21").Value = '  +0.01%  '

# Row 22 - Uniswap
$ws.Range("E22").Value = '  -0.46%  '

# Row 23 - Avalanche
$ws.Range("E23").Value = '  +4.29%  '

# Row 24 - Toncoin
$ws.Range("E24").Value = '  -2.77%  '

# Row 25 - Monero
Set-TextValue "D25" '149.75'
$ws.Range("E25").Value = '  +2.03%  '

# Row 26 - Cosmos
$ws.Range("E26").Value = '  -0.14%  '

# Row 27 - Stellar
$ws.Range("E27").Value = '  -0.91%  '

# Row 28 - BinanceUSD
$ws.Range("E28").Value = '  +0.00%  '

# Row 29 - Ethereum Classic
Set-TextValue "D29" '15.63'
$ws.Range("E29").Value = '  +0.82%  '

# Row 30 - PancakeSwap
$ws.Range("E30").Value = '  +0.34%  '

# Row 31 - Hedera
$ws.Range("E31").Value = '  +0.82%  '

# Row 32 - Filecoin
$ws.Range("E32").Value = '  +0.76%  '

# Row 33 - Maker
Set-TextValue "D33" '1.445.27'
$ws.Range("E33").Value = '  +2.53%  '

# Row 34 - Internet Computer (DFINITY)
$ws.Range("E34").Value = '  -0.03%  '

# Row 35 - Lido DAO Token
$ws.Range("E35").Value = '  -0.42%  '

# Row 36 - Huobi Token
Set-TextValue "D36" '2.37'
$ws.Range("E36").Value = '  -0.18%  '

# Row 37 - Immutable X
$ws.Range("E37").Value = '  +0.94%  '

# Row 38 - ARBITRUM
Set-TextValue "D38" '0.876'
$ws.Range("E38").Value = '  -0.11%  '

# Row 39 - VeChain
$ws.Range("E39").Value = '  +0.52%  '

# Row 40 - Trust Wallet Token
Set-TextValue "D40" '0.913'
$ws.Range("E40").Value = '  +15.08%  '

# Row 41 - Aave
Set-TextValue "D41" '70.15'
$ws.Range("E41").Value = '  +9.35%  '

# Row 42 - WEMIX Token
$ws.Range("E42").Value = '  -0.11%  '

# Row 43 - Pax Dollar
$ws.Range("E43").Value = '  +0.06%  '

# Row 45 - mCoin
$ws.Range("E45").Value = '  +0.57%  '

# Row 46 - MX Token
$ws.Range("E46").Value = '  +0.19%  '

# Row 47 - Rocket Pool ETH
Set-TextValue "D47" '1.781.65'
$ws.Range("E47").Value = '  +0.08%  '

# Row 48 - Render Token
$ws.Range("E48").Value = '  +4.04%  '

# Row 49 - Quant
Set-TextValue "D49" '85.96'
$ws.Range("E49").Value = '  -1.96%  '

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = '  -0.54%  '

# Row 51 - Algorand
Set-TextValue "D51" '0.0990'
$ws.Range("E51").Value = '  +0.69%  '
